# Case_0_67 / res_bus / vm_pu.xlsx — update voltage results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043141713594257
$ws.Range("D2").Value = 1.052482017688416
$ws.Range("E2").Value = 1.041156550327732
$ws.Range("F2").Value = 1.059408031509647
$ws.Range("I2").Value = 1.034717511710976
$ws.Range("J2").Value = 1.048213753378079
$ws.Range("K2").Value = 1.055230628717997
$ws.Range("L2").Value = 1.043936862582978
$ws.Range("M2").Value = 1.062137617646028
$ws.Range("N2").Value = 1.019843307565464
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044482344636537
$ws.Range("D3").Value = 1.05379617344803
$ws.Range("E3").Value = 1.042308310228442
$ws.Range("F3").Value = 1.060872610398774
$ws.Range("I3").Value = 1.034933925769033
$ws.Range("J3").Value = 1.049199170057318
$ws.Range("K3").Value = 1.056356352825202
$ws.Range("L3").Value = 1.044898296560397
$ws.Range("M3").Value = 1.063414775286938
$ws.Range("N3").Value = 1.020179458231152
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045348882553337
$ws.Range("D4").Value = 1.054645919465514
$ws.Range("E4").Value = 1.043053050402776
$ws.Range("F4").Value = 1.061819947172719
$ws.Range("I4").Value = 1.035072080063426
$ws.Range("J4").Value = 1.049835407551656
$ws.Range("K4").Value = 1.057083627388556
$ws.Range("L4").Value = 1.045519308542916
$ws.Range("M4").Value = 1.064240322906425
$ws.Range("N4").Value = 1.020396274279986
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045712954610027
$ws.Range("D5").Value = 1.055003012834022
$ws.Range("E5").Value = 1.043366016441585
$ws.Range("F5").Value = 1.062218129865704
$ws.Range("I5").Value = 1.035129710796643
$ws.Range("J5").Value = 1.050102551425306
$ws.Range("K5").Value = 1.057389103224557
$ws.Range("L5").Value = 1.04578012147027
$ws.Range("M5").Value = 1.064587181586768
$ws.Range("N5").Value = 1.020487258091707
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045774071113191
$ws.Range("D6").Value = 1.055062962294645
$ws.Range("E6").Value = 1.043418557708707
$ws.Range("F6").Value = 1.062284982069589
$ws.Range("I6").Value = 1.035139360923274
$ws.Range("J6").Value = 1.050147386774437
$ws.Range("K6").Value = 1.057440378208194
$ws.Range("L6").Value = 1.045823897914622
$ws.Range("M6").Value = 1.064645408984596
$ws.Range("N6").Value = 1.020502524987343
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045353748162674
$ws.Range("D7").Value = 1.054650691509245
$ws.Range("E7").Value = 1.043057232748744
$ws.Range("F7").Value = 1.061825268008903
$ws.Range("I7").Value = 1.035072851892737
$ws.Range("J7").Value = 1.049838978435444
$ws.Range("K7").Value = 1.057087710229701
$ws.Range("L7").Value = 1.045522794558115
$ws.Range("M7").Value = 1.064244958435976
$ws.Range("N7").Value = 1.020397490659873
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043594983516239
$ws.Range("D8").Value = 1.052926269519644
$ws.Range("E8").Value = 1.041545903641882
$ws.Range("F8").Value = 1.059903066739545
$ws.Range("I8").Value = 1.034791039596699
$ws.Range("J8").Value = 1.048547069866396
$ws.Range("K8").Value = 1.055611311499693
$ws.Range("L8").Value = 1.044262013200143
$ws.Range("M8").Value = 1.062569419310556
$ws.Range("N8").Value = 1.019957055921425
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040488411403743
$ws.Range("D9").Value = 1.049882840552151
$ws.Range("E9").Value = 1.038878585530444
$ws.Range("F9").Value = 1.056513063723964
$ws.Range("I9").Value = 1.034280016693602
$ws.Range("J9").Value = 1.046259753032331
$ws.Range("K9").Value = 1.053000785999842
$ws.Range("L9").Value = 1.042031801497706
$ws.Range("M9").Value = 1.059610145901546
$ws.Range("N9").Value = 1.019175581047047
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038412096900294
$ws.Range("D10").Value = 1.047850426266706
$ws.Range("E10").Value = 1.037097378375161
$ws.Range("F10").Value = 1.054250862203051
$ws.Range("I10").Value = 1.033929580801075
$ws.Range("J10").Value = 1.044727409574899
$ws.Range("K10").Value = 1.051254212196257
$ws.Range("L10").Value = 1.040539060722382
$ws.Range("M10").Value = 1.057632495339247
$ws.Range("N10").Value = 1.018650924752226
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.03751171587252
$ws.Range("D11").Value = 1.046969491085067
$ws.Range("E11").Value = 1.036325341250355
$ws.Range("F11").Value = 1.053270716879101
$ws.Range("I11").Value = 1.033775513221056
$ws.Range("J11").Value = 1.044062072400494
$ws.Range("K11").Value = 1.050496399657132
$ws.Range("L11").Value = 1.039891239130411
$ws.Range("M11").Value = 1.056774954063185
$ws.Range("N11").Value = 1.018422857733476
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037177069810005
$ws.Range("D12").Value = 1.046642134423112
$ws.Range("E12").Value = 1.036038453790699
$ws.Range("F12").Value = 1.052906551830964
$ws.Range("I12").Value = 1.03371793507485
$ws.Range("J12").Value = 1.043814658809328
$ws.Range("K12").Value = 1.050214679329385
$ws.Range("L12").Value = 1.039650387119533
$ws.Range("M12").Value = 1.056456238856841
$ws.Range("N12").Value = 1.018338008933224
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037248861870034
$ws.Range("D13").Value = 1.046712359867886
$ws.Range("E13").Value = 1.036099997550004
$ws.Range("F13").Value = 1.052984670880778
$ws.Range("I13").Value = 1.033730301664475
$ws.Range("J13").Value = 1.043867742496599
$ws.Range("K13").Value = 1.05027512000257
$ws.Range("L13").Value = 1.039702060808297
$ws.Range("M13").Value = 1.056524612858041
$ws.Range("N13").Value = 1.018356215397741
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037484058128432
$ws.Range("D14").Value = 1.04694243452741
$ws.Range("E14").Value = 1.036301629483472
$ws.Range("F14").Value = 1.053240616862652
$ws.Range("I14").Value = 1.033770760951523
$ws.Range("J14").Value = 1.044041626804651
$ws.Range("K14").Value = 1.050473117386968
$ws.Range("L14").Value = 1.039871334782526
$ws.Range("M14").Value = 1.05674861279259
$ws.Range("N14").Value = 1.01841584685724
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037628943042613
$ws.Range("D15").Value = 1.047084172669019
$ws.Range("E15").Value = 1.036425845752899
$ws.Range("F15").Value = 1.053398300811747
$ws.Range("I15").Value = 1.033795642765747
$ws.Range("J15").Value = 1.044148725743725
$ws.Range("K15").Value = 1.050595078835362
$ws.Range("L15").Value = 1.039975600508438
$ws.Range("M15").Value = 1.056886601720235
$ws.Range("N15").Value = 1.018452569900568
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038471823863262
$ws.Range("D16").Value = 1.047908871937151
$ws.Range("E16").Value = 1.037148599534425
$ws.Range("F16").Value = 1.054315898151654
$ws.Range("I16").Value = 1.033939756652065
$ws.Range("J16").Value = 1.044771527080334
$ws.Range("K16").Value = 1.05130447300736
$ws.Range("L16").Value = 1.040582023543547
$ws.Range("M16").Value = 1.057689381697693
$ws.Range("N16").Value = 1.01866604200778
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039000182495534
$ws.Range("D17").Value = 1.048425943209378
$ws.Range("E17").Value = 1.037601756942639
$ws.Range("F17").Value = 1.05489131882116
$ws.Range("I17").Value = 1.034029531867869
$ws.Range("J17").Value = 1.045161703227555
$ws.Range("K17").Value = 1.051749043245402
$ws.Range("L17").Value = 1.040962024594049
$ws.Range("M17").Value = 1.058192618038308
$ws.Range("N17").Value = 1.018799708974143
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039308237867872
$ws.Range("D18").Value = 1.048727456947504
$ws.Range("E18").Value = 1.03786600267189
$ws.Range("F18").Value = 1.055226894736365
$ws.Range("I18").Value = 1.034081671828355
$ws.Range("J18").Value = 1.045389110812289
$ws.Range("K18").Value = 1.052008205784814
$ws.Range("L18").Value = 1.041183532640104
$ws.Range("M18").Value = 1.05848603131076
$ws.Range("N18").Value = 1.018877589119105
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039413255370998
$ws.Range("D19").Value = 1.048830251039896
$ws.Range("E19").Value = 1.037956091318084
$ws.Range("F19").Value = 1.055341307917028
$ws.Range("I19").Value = 1.034099412176203
$ws.Range("J19").Value = 1.045466621262682
$ws.Range("K19").Value = 1.052096548585036
$ws.Range("L19").Value = 1.041259037498941
$ws.Range("M19").Value = 1.058586058136243
$ws.Range("N19").Value = 1.018904129764806
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038943507822261
$ws.Range("D20").Value = 1.04837047520387
$ws.Range("E20").Value = 1.037553145054783
$ws.Range("F20").Value = 1.0548295875912
$ws.Range("I20").Value = 1.034019923056024
$ws.Range("J20").Value = 1.045119859208417
$ws.Range("K20").Value = 1.051701360383216
$ws.Range("L20").Value = 1.04092126860012
$ws.Range("M20").Value = 1.058138637575156
$ws.Range("N20").Value = 1.018785376630311
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037414804382241
$ws.Range("D21").Value = 1.046874687127055
$ws.Range("E21").Value = 1.03624225720726
$ws.Range("F21").Value = 1.053165249830832
$ws.Range("I21").Value = 1.033758856386332
$ws.Range("J21").Value = 1.043990429908761
$ws.Range("K21").Value = 1.050414818615576
$ws.Range("L21").Value = 1.039821493980015
$ws.Range("M21").Value = 1.056682655575722
$ws.Range("N21").Value = 1.018398290603038
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036452462084263
$ws.Range("D22").Value = 1.045933424252941
$ws.Range("E22").Value = 1.035417362446581
$ws.Range("F22").Value = 1.052118259615336
$ws.Range("I22").Value = 1.033592684492273
$ws.Range("J22").Value = 1.043278703983863
$ws.Range("K22").Value = 1.049604556824554
$ws.Range("L22").Value = 1.039128734471031
$ws.Range("M22").Value = 1.055766142105212
$ws.Range("N22").Value = 1.018154135142616
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036962732243371
$ws.Range("D23").Value = 1.046432483146271
$ws.Range("E23").Value = 1.035854721250747
$ws.Range("F23").Value = 1.052673343464617
$ws.Range("I23").Value = 1.033680968029004
$ws.Range("J23").Value = 1.04365615719399
$ws.Range("K23").Value = 1.05003422252036
$ws.Range("L23").Value = 1.039496102658042
$ws.Range("M23").Value = 1.056252107291779
$ws.Range("N23").Value = 1.01828364083748
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.0389691170605
$ws.Range("D24").Value = 1.048395539070907
$ws.Range("E24").Value = 1.037575110900792
$ws.Range("F24").Value = 1.054857481453312
$ws.Range("I24").Value = 1.034024265557821
$ws.Range("J24").Value = 1.045138767261742
$ws.Range("K24").Value = 1.05172290667245
$ws.Range("L24").Value = 1.040939684912265
$ws.Range("M24").Value = 1.058163029380927
$ws.Range("N24").Value = 1.018791853063794
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041292441770002
$ws.Range("D25").Value = 1.050670231626194
$ws.Range("E25").Value = 1.039568664649417
$ws.Range("F25").Value = 1.057389827020154
$ws.Range("I25").Value = 1.03441384389796
$ws.Range("J25").Value = 1.046852380502701
$ws.Range("K25").Value = 1.053676749793877
$ws.Range("L25").Value = 1.042609397003005
$ws.Range("M25").Value = 1.060376016400234
$ws.Range("N25").Value = 1.019937825394386
